$d = $word.ActiveDocument

# --- Paragraph: contact form ---
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()
$paraStart = $d.Content.End - 1
$newPara = $d.Range($paraStart, $paraStart)
$url = "https://www.w3schools.com/howto/howto_css_contact_form.asp"
$newPara.InsertAfter($url + " contact form")
$linkRange = $d.Range($paraStart, $paraStart + $url.Length)
$link = $d.Hyperlinks.Add($linkRange, $url)

# --- Paragraph: skill box ---
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()
$paraStart = $d.Content.End - 1
$newPara = $d.Range($paraStart, $paraStart)
$url = "https://www.w3schools.com/howto/howto_css_skill_bar.asp"
$newPara.InsertAfter($url + " skill box")
$linkRange = $d.Range($paraStart, $paraStart + $url.Length)
$link = $d.Hyperlinks.Add($linkRange, $url)
